$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ref, $val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "36.096.04"
Set-TextCell "E2" "  -1.24%  "
Set-TextCell "D3" "2.009.45"
Set-TextCell "E3" "  -1.60%  "
Set-TextCell "E4" "  +0.20%  "
Set-TextCell "D5" "251.32"
Set-TextCell "E5" "  +2.79%  "
Set-TextCell "E6" "  -2.96%  "
Set-TextCell "D7" "62.46"
Set-TextCell "E7" "  +15.82%  "
Set-TextCell "E8" "  +0.31%  "
Set-TextCell "D9" "59.10"
Set-TextCell "E9" "  -5.98%  "
Set-TextCell "E10" "  +1.39%  "
Set-TextCell "D11" "0.0746"
Set-TextCell "E11" "  +0.21%  "
Set-TextCell "E12" "  -1.36%  "
Set-TextCell "D13" "0.923"
Set-TextCell "E13" "  -1.63%  "
Set-TextCell "D14" "14.87"
Set-TextCell "E14" "  +3.24%  "
Set-TextCell "D15" "2.308.65"
Set-TextCell "E15" "  -1.18%  "
Set-TextCell "D16" "5.41"
Set-TextCell "E16" "  +0.12%  "
Set-TextCell "D17" "19.51"
Set-TextCell "E17" "  +14.37%  "
Set-TextCell "D18" "2.039.35"
Set-TextCell "E18" "  -0.09%  "
Set-TextCell "D19" "36.064.81"
Set-TextCell "E19" "  -0.99%  "
Set-TextCell "D20" "71.97"
Set-TextCell "E20" "  +1.14%  "
Set-TextCell "E21" "  +0.43%  "
Set-TextCell "D22" "5.25"
Set-TextCell "E22" "  +1.52%  "
Set-TextCell "D23" "233.73"
Set-TextCell "E23" "  -1.72%  "
Set-TextCell "E24" "  +22.73%  "
Set-TextCell "E25" "  -0.20%  "
Set-TextCell "E26" "  -2.65%  "
Set-TextCell "D27" "9.54"
Set-TextCell "E27" "  +4.01%  "
Set-TextCell "D28" "164.64"
Set-TextCell "E28" "  +0.51%  "
Set-TextCell "D29" "19.57"
Set-TextCell "E29" "  -2.05%  "
Set-TextCell "E30" "  -0.34%  "
Set-TextCell "D31" "5.15"
Set-TextCell "E31" "  +3.47%  "
Set-TextCell "D32" "1.19"
Set-TextCell "E32" "  +2.44%  "
Set-TextCell "E33" "  +24.17%  "
Set-TextCell "E34" "  +1.02%  "
Set-TextCell "D35" "4.47"
Set-TextCell "E35" "  +1.44%  "
Set-TextCell "D36" "2.48"
Set-TextCell "E36" "  +12.49%  "
Set-TextCell "E37" "  +0.10%  "
Set-TextCell "E38" "  -0.42%  "
Set-TextCell "D39" "5.74"
Set-TextCell "E39" "  +16.85%  "
Set-TextCell "D40" "0.111"
Set-TextCell "E40" "  +23.88%  "
Set-TextCell "E41" "  -1.13%  "
Set-TextCell "E42" "  +0.66%  "
Set-TextCell "D43" "0.0215"
Set-TextCell "E43" "  +1.20%  "
Set-TextCell "E44" "  +2.18%  "
Set-TextCell "D45" "16.81"
Set-TextCell "E45" "  +7.14%  "
Set-TextCell "B46" "FraxShare"
Set-TextCell "C46" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell "D46" "7.82"
Set-TextCell "E46" "  +5.18%  "
Set-TextCell "B47" "Aave"
Set-TextCell "C47" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell "D47" "93.77"
Set-TextCell "E47" "  +0.30%  "
Set-TextCell "D48" "1.416.95"
Set-TextCell "E48" "  +3.16%  "
Set-TextCell "D49" "2.35"
Set-TextCell "E49" "  +4.28%  "
Set-TextCell "D50" "2.90"
Set-TextCell "E50" "  -0.55%  "
Set-TextCell "D51" "47.04"
Set-TextCell "E51" "  +3.82%  "
